{"js": "// Update the date line and the 25 division-problem cells in the table.\n// A couple of \"before\" strings repeat (e.g. \"46\u00f72=\" appears twice in the\n// same row), so instead of a naive global find/replace we walk the\n// document's paragraphs in order and rewrite each one according to its\n// position \u2014 this is robust to duplicate values.\n\nconst replacements = [\n  \"2024-07-25 Thursday\", // 2024-07-24 Wednesday\n  \"13\u00f72=\", // 14\u00f76=\n  \"92\u00f72=\", // 40\u00f76=\n  \"40\u00f72=\", // 29\u00f72=\n  \"83\u00f78=\", // 46\u00f72=\n  \"55\u00f74=\", // 46\u00f72=\n  \"89\u00f79=\", // 81\u00f74=\n  \"73\u00f77=\", // 83\u00f73=\n  \"81\u00f73=\", // 41\u00f72=\n  \"94\u00f79=\", // 17\u00f75=\n  \"23\u00f72=\", // 31\u00f76=\n  \"76\u00f72=\", // 50\u00f79=\n  \"41\u00f73=\", // 62\u00f73=\n  \"73\u00f77=\", // 38\u00f78=\n  \"72\u00f78=\", // 49\u00f76=\n  \"56\u00f76=\", // 78\u00f74=\n  \"60\u00f75=\", // 33\u00f76=\n  \"98\u00f72=\", // 14\u00f78=\n  \"33\u00f73=\", // 79\u00f72=\n  \"91\u00f76=\", // 80\u00f75=\n  \"68\u00f78=\", // 86\u00f73=\n  \"52\u00f73=\", // 13\u00f79=\n  \"49\u00f75=\", // 10\u00f73=\n  \"33\u00f73=\", // 43\u00f72=\n  \"79\u00f78=\", // 45\u00f73=\n  \"16\u00f72=\", // 18\u00f76=\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Only paragraphs that actually contain text are meaningful here (the\n// empty table cells stay empty), and they occur in the same order as the\n// `replacements` array above.\nlet idx = 0;\nfor (let i = 0; i < paragraphs.items.length && idx < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.length === 0) {\n    continue;\n  }\n  para.insertText(replacements[idx], Word.InsertLocation.replace);\n  idx++;\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division-problem cells in the table.\n# A couple of \"before\" strings repeat (e.g. \"46\u00f72=\" appears twice in the\n# same row), so a straight Replace-All would collide; instead we search\n# occurrence-by-occurrence (left to right through the document) and feed\n# each match its own replacement text.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Occurrences {\n    param(\n        $Document,\n        [string]$OldText,\n        [string[]]$NewTexts\n    )\n\n    $rng = $Document.Range()\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $OldText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Wrap = 0  # wdFindStop - do not wrap back to start\n\n    $i = 0\n    while ($i -lt $NewTexts.Count -and $find.Execute()) {\n        $rng.Text = $NewTexts[$i]\n        $i++\n        $rng.Collapse(0)  # wdCollapseEnd - continue searching after the replacement\n    }\n}\n\nReplace-Occurrences $d \"2024-07-24 Wednesday\" @(\"2024-07-25 Thursday\")\n\nReplace-Occurrences $d \"14\u00f76=\" @(\"13\u00f72=\")\nReplace-Occurrences $d \"40\u00f76=\" @(\"92\u00f72=\")\nReplace-Occurrences $d \"29\u00f72=\" @(\"40\u00f72=\")\nReplace-Occurrences $d \"46\u00f72=\" @(\"83\u00f78=\", \"55\u00f74=\")\n\nReplace-Occurrences $d \"81\u00f74=\" @(\"89\u00f79=\")\nReplace-Occurrences $d \"83\u00f73=\" @(\"73\u00f77=\")\nReplace-Occurrences $d \"41\u00f72=\" @(\"81\u00f73=\")\nReplace-Occurrences $d \"17\u00f75=\" @(\"94\u00f79=\")\nReplace-Occurrences $d \"31\u00f76=\" @(\"23\u00f72=\")\n\nReplace-Occurrences $d \"50\u00f79=\" @(\"76\u00f72=\")\nReplace-Occurrences $d \"62\u00f73=\" @(\"41\u00f73=\")\nReplace-Occurrences $d \"38\u00f78=\" @(\"73\u00f77=\")\nReplace-Occurrences $d \"49\u00f76=\" @(\"72\u00f78=\")\nReplace-Occurrences $d \"78\u00f74=\" @(\"56\u00f76=\")\n\nReplace-Occurrences $d \"33\u00f76=\" @(\"60\u00f75=\")\nReplace-Occurrences $d \"14\u00f78=\" @(\"98\u00f72=\")\nReplace-Occurrences $d \"79\u00f72=\" @(\"33\u00f73=\")\nReplace-Occurrences $d \"80\u00f75=\" @(\"91\u00f76=\")\nReplace-Occurrences $d \"86\u00f73=\" @(\"68\u00f78=\")\n\nReplace-Occurrences $d \"13\u00f79=\" @(\"52\u00f73=\")\nReplace-Occurrences $d \"10\u00f73=\" @(\"49\u00f75=\")\nReplace-Occurrences $d \"43\u00f72=\" @(\"33\u00f73=\")\nReplace-Occurrences $d \"45\u00f73=\" @(\"79\u00f78=\")\nReplace-Occurrences $d \"18\u00f76=\" @(\"16\u00f72=\")\n"}
